$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet lists shortage items in rows 7..23 (serial numbers 1..17), with a
# totals row and a footer (timestamp / page / developer credit) right below.
# This edit inserts two new item rows:
#   - "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML" right after "KEMIFORGE ..."
#     (i.e. before "OPLEX-N SYRUP 125ML"), becoming the new serial #11.
#   - "OSTEO ZAD SYRUP 120 ML" right after "OPLEX-N SYRUP 125ML"
#     (i.e. before "PLEGICA 1% EYE DROPS 10 ML"), becoming the new serial #13.
# Every row after each insertion point shifts down by one, the running total
# grows by the two new sell prices, and the footer timestamp is refreshed.
# ---------------------------------------------------------------------------

# 1) Insert the two new rows at the right spots. Doing the earlier (row 17)
#    insert first means the second insert position (row 19) already accounts
#    for the first row having been added.
$ws.Rows("17:17").Insert()
$ws.Rows("19:19").Insert()

# 2) The freshly inserted rows come back with generic/default formatting
#    (no borders, wrong number formats). Clone the look of a neighboring,
#    still-correctly-styled item row (row 18, "OPLEX-N ...") onto them so
#    they match the rest of the table exactly.
$ws.Range("A18:Q18").Copy()
$ws.Range("A17:Q17").PasteSpecial(-4122)
$ws.Range("A19:Q19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Row heights: odd-positioned item rows are 25.5pt, even-positioned ones
#    are 24.75pt in this sheet's alternating rhythm.
$ws.Rows("17:17").RowHeight = 25.5
$ws.Rows("19:19").RowHeight = 25.5

# 4) Recreate the merged cells for the new rows (serial+blank / name / stock /
#    limit / price), matching every other item row's merge layout.
$ws.Range("A17:B17").Merge()
$ws.Range("C17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()
$ws.Range("N17:O17").Merge()

$ws.Range("A19:B19").Merge()
$ws.Range("C19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("N19:O19").Merge()

# 5) Fill in the values for the two new items.
#    Row 17: MAXILASE 200 CEIP UNIT/ML SYRUP 100ML  -> serial 11
$ws.Range("A17").Value = 11
$ws.Range("C17").Value = "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML"
$ws.Range("H17").Value = "4:0"
$ws.Range("L17").Value = "1"
$ws.Range("N17").Value = "57.00"
$ws.Range("P17").Value = "57.0000"
$ws.Range("Q17").Value = "1:0"

#    Row 19: OSTEO ZAD SYRUP 120 ML -> serial 13
$ws.Range("A19").Value = 13
$ws.Range("C19").Value = "OSTEO ZAD SYRUP 120 ML"
$ws.Range("H19").Value = "1:0"
$ws.Range("L19").Value = "1"
$ws.Range("N19").Value = "70.00"
$ws.Range("P19").Value = "70.0000"
$ws.Range("Q19").Value = "1:0"

# 6) Renumber the serials of every item row that followed the insert points
#    (they were copied down unchanged by Insert, so just reassign 12..19).
$ws.Range("A18").Value = 12
$ws.Range("A20").Value = 14
$ws.Range("A21").Value = 15
$ws.Range("A22").Value = 16
$ws.Range("A23").Value = 17
$ws.Range("A24").Value = 18
$ws.Range("A25").Value = 19

# 7) The running total (now on row 26) grows by the two new sell prices.
$ws.Range("P26").Value = 1507

# 8) Refresh the generated-at timestamp in the footer (now row 27).
$ws.Range("A27").Value = "Monday, 22 September, 2025 11:55 AM"
